# Update 'F' column (想去人数 / want-to-go count) values across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 2891
$ws.Range("F6").Value = 2498
$ws.Range("F7").Value = 18
$ws.Range("F8").Value = 39
$ws.Range("F10").Value = 2912
$ws.Range("F11").Value = 348
$ws.Range("F13").Value = 7223
$ws.Range("F14").Value = 316
$ws.Range("F15").Value = 41
$ws.Range("F16").Value = 98
$ws.Range("F17").Value = 225
$ws.Range("F18").Value = 115
$ws.Range("F19").Value = 479
$ws.Range("F20").Value = 8611
$ws.Range("F23").Value = 263
$ws.Range("F24").Value = 66
$ws.Range("F28").Value = 93
$ws.Range("F31").Value = 40
$ws.Range("F33").Value = 95
$ws.Range("F34").Value = 2605
$ws.Range("F37").Value = 35
$ws.Range("F38").Value = 1175
$ws.Range("F40").Value = 706
$ws.Range("F41").Value = 3714
$ws.Range("F42").Value = 7
$ws.Range("F43").Value = 185
$ws.Range("F44").Value = 22
$ws.Range("F45").Value = 1193
$ws.Range("F46").Value = 195
$ws.Range("F47").Value = 35
$ws.Range("F48").Value = 9
$ws.Range("F49").Value = 27

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 29
$ws.Range("F5").Value = 250

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 2891
$ws.Range("F4").Value = 29
$ws.Range("F5").Value = 250
$ws.Range("F7").Value = 2498
$ws.Range("F9").Value = 18
$ws.Range("F10").Value = 39
$ws.Range("F12").Value = 2912
$ws.Range("F13").Value = 348
$ws.Range("F17").Value = 7223
$ws.Range("F18").Value = 316
$ws.Range("F19").Value = 98
$ws.Range("F20").Value = 225
$ws.Range("F21").Value = 115
$ws.Range("F22").Value = 479
$ws.Range("F23").Value = 8612
$ws.Range("F25").Value = 263
$ws.Range("F26").Value = 66
$ws.Range("F29").Value = 93
$ws.Range("F31").Value = 40
$ws.Range("F34").Value = 95
$ws.Range("F35").Value = 2605
$ws.Range("F38").Value = 35
$ws.Range("F39").Value = 1175
$ws.Range("F40").Value = 706
$ws.Range("F42").Value = 3714
$ws.Range("F43").Value = 185
$ws.Range("F44").Value = 22
$ws.Range("F46").Value = 1193
$ws.Range("F47").Value = 195
$ws.Range("F48").Value = 35
$ws.Range("F49").Value = 27
